# Limpieza de documento: eliminar separadores decorativos, imagenes
# inline y los parrafos vacios de espaciado ("before=40 twips") que
# Word insertaba justo antes de cada encabezado Heading3 tras una tabla.
#
# Estrategia: recorrer la coleccion Paragraphs de atras hacia adelante
# (para que borrar un parrafo no desplace los indices de los que aun
# quedan por visitar) y borrar el rango completo de cada parrafo que
# cumpla una de estas tres condiciones:
#   1. Contiene una imagen en linea (InlineShapes.Count > 0)
#   2. Su texto es la linea separadora de guiones Unicode (U+2500 "-")
#   3. Esta vacio (solo la marca de parrafo) y tiene
#      Format.SpaceBefore = 2pt (= 40 twips)

$d = $word.ActiveDocument

$separatorChar = [char]0x2500

$deleted = 0
$count = $d.Paragraphs.Count

for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $rng = $p.Range
    $text = $rng.Text

    $hasDrawing = ($rng.InlineShapes.Count -gt 0)
    $isSeparator = ($text.IndexOf($separatorChar) -ge 0)
    $isEmpty = ($text.Length -le 1)
    $spaceBefore = $p.Format.SpaceBefore
    $isEmptySpacer = ($isEmpty -and $spaceBefore -eq 2)

    if ($hasDrawing -or $isSeparator -or $isEmptySpacer) {
        $rng.Delete()
        $deleted = $deleted + 1
    }
}

Write-Output "Parrafos eliminados: $deleted"
Write-Output "Parrafos restantes: $($d.Paragraphs.Count)"
